$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$shape = $s.Shapes.Item(2)
$para = $shape.TextFrame.TextRange.Paragraphs(1)
$para.Runs(1).Text = "Fetch raw dataset (from Internet or disk file)"
